$wb = $excel.ActiveWorkbook

$ws = $wb.Sheets.Item("ALC")
$ws.Range("H17").Value = 340.1
$ws.Range("J17").Value = 340.1
$ws.Range("L17").Value = 1020.3
$ws.Range("N17").Value = -1356.3
$ws.Range("H21").Value = 11413.889
$ws.Range("I21").Value = 8750
$ws.Range("J21").Value = 12175
$ws.Range("K21").Value = 8750
$ws.Range("L21").Value = 12175
$ws.Range("M21").Value = -8282
$ws.Range("N21").Value = -13111
$ws.Range("H23").Value = 11413.889
$ws.Range("I23").Value = 8750
$ws.Range("J23").Value = 12175
$ws.Range("K23").Value = 8750
$ws.Range("L23").Value = 12175
$ws.Range("M23").Value = -8516
$ws.Range("N23").Value = -12643
$ws.Range("H44").Value = 0
$ws.Range("J44").Value = 0
$ws.Range("L44").Value = 0
$ws.Range("N44").ClearContents()
$ws.Range("H61").Value = 7185765.5
$ws.Range("I61").Value = 28571428
$ws.Range("J61").Value = 57211.332
$ws.Range("K61").Value = 85714284
$ws.Range("L61").Value = 171633.996
$ws.Range("M61").Value = -85714112
$ws.Range("N61").Value = -171977.996
$ws.Range("H100").Value = 1108.3572
$ws.Range("I100").Value = 1191.5834
$ws.Range("J100").Value = 609
$ws.Range("K100").Value = 1191.5834
$ws.Range("L100").Value = 609
$ws.Range("M100").Value = -650.5834
$ws.Range("N100").Value = -1691
$ws.Range("H129").Value = 350621.75
$ws.Range("J129").Value = 376312.7
$ws.Range("L129").Value = 1128938.1
$ws.Range("N129").Value = -1138938.1
$ws.Range("H137").Value = 2475.2942
$ws.Range("I137").Value = 1659.7
$ws.Range("K137").Value = 4979.1
$ws.Range("M137").Value = -2429.1
$ws.Range("H138").Value = 2325.303
$ws.Range("I138").Value = 1465.96
$ws.Range("J138").Value = 2615.6216
$ws.Range("K138").Value = 4397.88
$ws.Range("L138").Value = 7846.864799999999
$ws.Range("M138").Value = 742.1199999999999
$ws.Range("N138").Value = -18126.8648

$ws = $wb.Sheets.Item("ARM")
$ws.Range("H2").Value = 49184.19
$ws.Range("I2").Value = 1558.6111
$ws.Range("K2").Value = 1558.6111
$ws.Range("M2").Value = -1445.6111
$ws.Range("H59").Value = 0
$ws.Range("I59").Value = 0
$ws.Range("J59").Value = 0
$ws.Range("K59").Value = 0
$ws.Range("L59").Value = 0
$ws.Range("M59").Value = 0
$ws.Range("N59").ClearContents()
$ws.Range("H116").Value = 49184.19
$ws.Range("I116").Value = 1558.6111
$ws.Range("K116").Value = 1558.6111
$ws.Range("M116").Value = 735.3888999999999

$ws = $wb.Sheets.Item("BSM")
$ws.Range("H3").Value = 49184.19
$ws.Range("I3").Value = 1558.6111
$ws.Range("K3").Value = 1558.6111
$ws.Range("M3").Value = -1444.6111
$ws.Range("H86").Value = 60034.26
$ws.Range("I86").Value = 74970.266
$ws.Range("J86").Value = 4024.25
$ws.Range("K86").Value = 74970.266
$ws.Range("L86").Value = 4024.25
$ws.Range("M86").Value = -73847.266
$ws.Range("N86").Value = -6270.25
$ws.Range("H89").Value = 60034.26
$ws.Range("I89").Value = 74970.266
$ws.Range("J89").Value = 4024.25
$ws.Range("K89").Value = 374851.33
$ws.Range("L89").Value = 20121.25
$ws.Range("M89").Value = -369235.33
$ws.Range("N89").Value = -31353.25
$ws.Range("H107").Value = 90951620
$ws.Range("I107").Value = 166743420
$ws.Range("K107").Value = 166743420
$ws.Range("M107").Value = -166741500

$ws = $wb.Sheets.Item("CRP")
$ws.Range("H86").Value = 1901.84
$ws.Range("I86").Value = 1519.3334
$ws.Range("J86").Value = 2475.6
$ws.Range("K86").Value = 1519.3334
$ws.Range("L86").Value = 2475.6
$ws.Range("M86").Value = -396.3334
$ws.Range("N86").Value = -4721.6
$ws.Range("H89").Value = 1901.84
$ws.Range("I89").Value = 1519.3334
$ws.Range("J89").Value = 2475.6
$ws.Range("K89").Value = 7596.666999999999
$ws.Range("L89").Value = 12378
$ws.Range("M89").Value = -1980.666999999999
$ws.Range("N89").Value = -23610
$ws.Range("H111").Value = 0
$ws.Range("J111").Value = 0
$ws.Range("L111").Value = 0
$ws.Range("N111").ClearContents()

$ws = $wb.Sheets.Item("CUL")
$ws.Range("H17").Value = 980.2
$ws.Range("I17").Value = 474.75
$ws.Range("K17").Value = 1424.25
$ws.Range("M17").Value = -1255.25
$ws.Range("H55").Value = 8195.967
$ws.Range("J55").Value = 5759.56
$ws.Range("L55").Value = 17278.68
$ws.Range("N55").Value = -17632.68
$ws.Range("H124").Value = 3585.4443
$ws.Range("I124").Value = 2941.8
$ws.Range("J124").Value = 4390
$ws.Range("K124").Value = 8825.400000000001
$ws.Range("L124").Value = 13170
$ws.Range("M124").Value = -3915.400000000001
$ws.Range("N124").Value = -22990
$ws.Range("H131").Value = 855.28
$ws.Range("I131").Value = 640
$ws.Range("J131").Value = 864.25
$ws.Range("K131").Value = 1920
$ws.Range("L131").Value = 2592.75
$ws.Range("M131").Value = 3120
$ws.Range("N131").Value = -12672.75

$ws = $wb.Sheets.Item("GSM")
$ws.Range("H5").Value = 1254250
$ws.Range("J5").Value = 5666.6665
$ws.Range("L5").Value = 5666.6665
$ws.Range("N5").Value = -5890.6665
$ws.Range("H80").Value = 100002040
$ws.Range("J80").Value = 1986.5714
$ws.Range("L80").Value = 1986.5714
$ws.Range("N80").Value = -3982.5714
$ws.Range("H83").Value = 100002040
$ws.Range("J83").Value = 1986.5714
$ws.Range("L83").Value = 9932.857
$ws.Range("N83").Value = -19916.857
$ws.Range("H120").Value = 34295.75
$ws.Range("J120").Value = 34295.75
$ws.Range("L120").Value = 34295.75
$ws.Range("N120").Value = -43971.75

$ws = $wb.Sheets.Item("LTW")
$ws.Range("H61").Value = 2708.1
$ws.Range("I61").Value = 2098.75
$ws.Range("J61").Value = 3622.125
$ws.Range("K61").Value = 2098.75
$ws.Range("L61").Value = 3622.125
$ws.Range("M61").Value = -1896.75
$ws.Range("N61").Value = -4026.125
$ws.Range("H113").Value = 2708.1
$ws.Range("I113").Value = 2098.75
$ws.Range("J113").Value = 3622.125
$ws.Range("K113").Value = 2098.75
$ws.Range("L113").Value = 3622.125
$ws.Range("M113").Value = 71.25
$ws.Range("N113").Value = -7962.125

$ws = $wb.Sheets.Item("WVR")
$ws.Range("H2").Value = 852304.2
$ws.Range("I2").Value = 1441357.1
$ws.Range("J2").Value = 27630
$ws.Range("K2").Value = 1441357.1
$ws.Range("L2").Value = 27630
$ws.Range("M2").Value = -1441245.1
$ws.Range("N2").Value = -27854
$ws.Range("H28").Value = 11631.143
$ws.Range("J28").Value = 11631.143
$ws.Range("L28").Value = 11631.143
$ws.Range("N28").Value = -12327.143
